# Generate Report for Handoff
#
# Re-runs the localization-status handoff report generation: the four
# files that were sitting at "Ready for handoff" (rows 4-7 on each
# table) just had a fresh handoff xliff generated for them, so their
# Priority flips from "low" to "ht" and the handoff timestamps move
# forward a few seconds.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newGenerateDate = "2016-08-23 10:32:20"
$newZhCnHandoff  = "2016-08-23 10:32:13"
$newDeDeHandoff  = "2016-08-23 10:32:20"

foreach ($row in 4..7) {
    # Overview: "Latest HO Xliff Generate Date"
    $overview.Range("G$row").Value = $newGenerateDate

    # zh-cn table: Priority + Latest Handoff Datetime
    $zhcn.Range("E$row").Value = "ht"
    $zhcn.Range("H$row").Value = $newZhCnHandoff

    # de-de table: Priority + Latest Handoff Datetime
    $dede.Range("E$row").Value = "ht"
    $dede.Range("H$row").Value = $newDeDeHandoff
}
